$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.499.20"
$ws.Range("E2").Value = "  +5.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.181.94"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "403.89"
$ws.Range("E5").Value = "  +4.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.20"
$ws.Range("E6").Value = "  +5.24%  "
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  +5.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.97"
$ws.Range("E10").Value = "  +4.97%  "
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0883"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.672.03"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.23"
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.06"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("E16").Value = "  +8.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.187.33"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.49"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "54.329.29"
$ws.Range("E19").Value = "  +4.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.31"
$ws.Range("E20").Value = "  +3.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.88"
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0996"
$ws.Range("E22").Value = "  +2.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.63"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "274.27"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.29"
$ws.Range("E25").Value = "  +4.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.01"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.75"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.170"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +3.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.08"
$ws.Range("E32").Value = "  +6.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0499"
$ws.Range("E33").Value = "  +10.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "36.85"
$ws.Range("E34").Value = "  +3.28%  "
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.68"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("E37").Value = "  +8.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.997"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.84"
$ws.Range("E39").Value = "  +9.46%  "
$ws.Range("E40").Value = "  +10.52%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.93"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.291"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.30"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.23"
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.21"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.087.40"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0343"
$ws.Range("E50").Value = "  +7.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0506"
$ws.Range("E51").Value = "  +8.62%  "
